# dev-1.0.0 : Update Templete Sales dan Labor
#
# Clears the monthly (Jan-Dec, columns E:P) data values that were
# previously filled in the "Sales" worksheet, leaving only E2 = 1.
# The SUM formulas in Q/R/S recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sales")

# Row 2 keeps its value in column E, but F2:P2 get cleared.
$ws.Range("F2:P2").ClearContents()

# Rows 3-22 get their whole E:P block cleared.
$ws.Range("E3:P22").ClearContents()

# Update the selection/view so only E2 is selected (no multi-cell
# range selected, and no horizontal scroll to column I).
$ws.Range("E2").Select()
